$wb = $excel.ActiveWorkbook

# --- Sheet "BrokerNewPaymentData" --------------------------------------
# LoadID/Invoice values move from the 08/06/2018 batch to the 08/07/2018
# batch, and the second row's Amount is bumped from 100 to 200.
$wsNew = $wb.Worksheets.Item("BrokerNewPaymentData")
$wsNew.Range("B2").Value = "080718A07"
$wsNew.Range("C2").Value = "080718A07"
$wsNew.Range("B3").Value = "080718A08"
$wsNew.Range("C3").Value = "080718A08"
$wsNew.Range("D3").Value = 200

# --- Sheet "BrokerPaymentDataforUnmatchedCr" ----------------------------
# Same 08/06/2018 -> 08/07/2018 date roll for the unmatched-credit email
# and invoice/loadid values.
$wsUnmatched = $wb.Worksheets.Item("BrokerPaymentDataforUnmatchedCr")
$wsUnmatched.Range("A2").Value = "umCVK080718A01@loadpaytest.truckstop.com"
$wsUnmatched.Range("B2").Value = "080718A01UM"
$wsUnmatched.Range("C2").Value = "080718A01UM"

# --- Selections (smart-wait landed on different cells after the edits) -
$wsNew.Range("C3").Select()
$wsUnmatched.Range("A2").Select()

# --- Active sheet/tab ----------------------------------------------------
$wsUpdated = $wb.Worksheets.Item("BrokerUpdatedPaymentData")
$wsUpdated.Activate()
